$d = $word.ActiveDocument

$d.Content.Find.Execute("obvezna je vratiti se na rad dana ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ obv }} je vratiti se na rad dana ", 2)

$d.Content.Find.Execute("NAPUTAK O ZAŠTITI PRAVA: Protiv ove Odluke radnik/ca", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NAPUTAK O ZAŠTITI PRAVA: Protiv ove Odluke {{ radn }}", 2)

$d.Content.Find.Execute("1. Radniku/ci:primio", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1. {{ Rdn }}: {{ prim }}", 2)
